$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster / Ligand / Receptor / Target cluster / then 16 numeric metric columns (E..T)
$data = @(
    @(2, "ECs", "Cd14", "Itga4", "ECs", 3, 1, 1.029662, 3.088986, 0.008285277389098726, 0.008285277389098728, 2, 0.6666666666666666, 0.3809099999999999, 1.14273, 0.006635732896411959, 0.006635732896411961, 0.39220855242, 3.52987697178, 0.0000549788877267406, 0.00005497888772674064),
    @(3, "ECs", "Cd14", "Itga4", "FAPs", 3, 1, 1.029662, 3.088986, 0.008285277389098726, 0.008285277389098728, 3, 1, 0.3194813333333333, 0.9584440000000001, 0.005565600255676026, 0.005565600255676028, 0.3289577886426667, 2.960620097784, 0.00004611254195511466, 0.00004611254195511469),
    @(4, "ECs", "Cd14", "Itga4", "M2", 3, 1, 1.029662, 3.088986, 0.008285277389098726, 0.008285277389098728, 3, 1, 54.69403966666666, 164.082119, 0.9528104755815301, 0.9528104755815303, 56.31637427125933, 506.847368441334, 0.007894299089432056, 0.007894299089432059),
    @(5, "ECs", "Cd14", "Itga4", "sCs", 3, 1, 1.029662, 3.088986, 0.008285277389098726, 0.008285277389098728, 3, 1, 2.008422, 6.025265999999999, 0.03498819126638183, 0.03498819126638184, 2.067995813364, 18.611962320276, 0.0002898868699848149, 0.000289886869984815),
    @(6, "FAPs", "Cd14", "Itga4", "ECs", 3, 1, 1.017189, 3.051567, 0.00818491215771772, 0.008184912157717722, 2, 0.6666666666666666, 0.3809099999999999, 1.14273, 0.006635732896411959, 0.006635732896411961, 0.3874574619899999, 3.487117157909999, 0.00005431289085920966, 0.0000543128908592097),
    @(7, "FAPs", "Cd14", "Itga4", "FAPs", 3, 1, 1.017189, 3.051567, 0.00818491215771772, 0.008184912157717722, 3, 1, 0.3194813333333333, 0.9584440000000001, 0.005565600255676026, 0.005565600255676028, 0.324972897972, 2.924756081748, 0.00004555394919767955, 0.00004555394919767958),
    @(8, "FAPs", "Cd14", "Itga4", "M2", 3, 1, 1.017189, 3.051567, 0.00818491215771772, 0.008184912157717722, 3, 1, 54.69403966666666, 164.082119, 0.9528104755815301, 0.9528104755815303, 55.63417551449699, 500.7075796304729, 0.007798670045588069, 0.007798670045588072),
    @(9, "FAPs", "Cd14", "Itga4", "sCs", 3, 1, 1.017189, 3.051567, 0.00818491215771772, 0.008184912157717722, 3, 1, 2.008422, 6.025265999999999, 0.03498819126638183, 0.03498819126638184, 2.042944765758, 18.386502891822, 0.0002863752720727616, 0.0002863752720727617),
    @(10, "M2", "Cd14", "Itga4", "ECs", 3, 1, 120.9972243333333, 362.991673, 0.9736161642487271, 0.9736161642487271, 2, 0.6666666666666666, 0.3809099999999999, 1.14273, 0.006635732896411959, 0.006635732896411961, 46.08905272080999, 414.8014744872899, 0.006460656809583707, 0.006460656809583709),
    @(11, "M2", "Cd14", "Itga4", "FAPs", 3, 1, 120.9972243333333, 362.991673, 0.9736161642487271, 0.9736161642487271, 3, 1, 0.3194813333333333, 0.9584440000000001, 0.005565600255676026, 0.005565600255676028, 38.65635455964578, 347.907191036812, 0.005418758372673027, 0.005418758372673029),
    @(12, "M2", "Cd14", "Itga4", "M2", 3, 1, 120.9972243333333, 362.991673, 0.9736161642487271, 0.9736161642487271, 3, 1, 54.69403966666666, 164.082119, 0.9528104755815301, 0.9528104755815303, 6617.826987243898, 59560.44288519508, 0.9276716804916948, 0.927671680491695),
    @(13, "M2", "Cd14", "Itga4", "sCs", 3, 1, 120.9972243333333, 362.991673, 0.9736161642487271, 0.9736161642487271, 3, 1, 2.008422, 6.025265999999999, 0.03498819126638183, 0.03498819126638184, 243.013487290002, 2187.121385610018, 0.03406506857477549, 0.03406506857477549),
    @(14, "sCs", "Cd14", "Itga4", "ECs", 3, 1, 1.232029333333333, 3.696088, 0.009913646204456457, 0.009913646204456457, 2, 0.6666666666666666, 0.3809099999999999, 1.14273, 0.006635732896411959, 0.006635732896411961, 0.4692922933599999, 4.22363064024, 0.00006578430824230126, 0.00006578430824230129),
    @(15, "sCs", "Cd14", "Itga4", "FAPs", 3, 1, 1.232029333333333, 3.696088, 0.009913646204456457, 0.009913646204456457, 3, 1, 0.3194813333333333, 0.9584440000000001, 0.005565600255676026, 0.005565600255676028, 0.3936103741191112, 3.542493367072, 0.00005517539185020452, 0.00005517539185020454),
    @(16, "sCs", "Cd14", "Itga4", "M2", 3, 1, 1.232029333333333, 3.696088, 0.009913646204456457, 0.009913646204456457, 3, 1, 54.69403966666666, 164.082119, 0.9528104755815301, 0.9528104755815303, 67.38466122783022, 606.461951050472, 0.009445825954815187, 0.00944582595481519),
    @(17, "sCs", "Cd14", "Itga4", "sCs", 3, 1, 1.232029333333333, 3.696088, 0.009913646204456457, 0.009913646204456457, 3, 1, 2.008422, 6.025265999999999, 0.03498819126638183, 0.03498819126638184, 2.474434817712, 22.269913359408, 0.0003468605495487628, 0.0003468605495487628)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}

Write-Output "Updated rows 2-17 (A1:T17)"